$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sponsor Work (row 13): add hours worked on Saturday (column G)
$ws.Range("G13").Value = 3
# Update Daily Total for Saturday column I13
$ws.Range("I13").Value = 3.5

# Daily Total row (row 14): Saturday total
$ws.Range("G14").Value = 3
# Weekly running total column I14
$ws.Range("I14").Value = 3.5
